$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test data for 2022 (was 2020/2021) - [IMP] New data for test environment
$ws.Range("D2").Value = "P1/2022/0001"
$ws.Range("D3").Value = 220123
$ws.Range("D4").Value = 22011214
$ws.Range("D5").Value = "IT/22/004"
$ws.Range("D8").Value = "P1/2022/0007"

# View tweaks: tab ratio, zoom and the active selection
$excel.ActiveWindow.TabRatio = 0.6
$excel.ActiveWindow.Zoom = 60
$ws.Range("A1").Select()
